$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Defect Log")

# --- Fill in newly-entered defect rows 76:81 -------------------------------
# Column map: B=Created Date, C=Title, D=Description, E=Status, H=Priority,
# I=Product Type, K=Product, M=Detected By, N=Assigned To, Q=Cause Analysis,
# S=Fixed Date.

$rowsData = @(
    @{ Row = 76; B = "22/10/2011"; C = "Section Manager";  D = " database: các item trong phần này đều dư"; N = "DungDV" },
    @{ Row = 77; B = "22/10/2011"; C = "Category Manager"; D = " database: các item trong phần này đều dư"; N = "DungDV" },
    @{ Row = 78; B = "22/10/2011"; C = "Front Page Manager"; D = " database: các item trong phần này đều dư"; N = "DungDV" },
    @{ Row = 79; B = "22/10/2011"; C = "MenuManager"; D = "Sheet Contents - Menu Manager:   chưa có link"; N = "ThiVT" },
    @{ Row = 80; B = "22/10/2011"; C = "MenuManager"; D = "Sheet Contents - Menu Trash:          chưa có link"; N = "ThiVT" },
    @{ Row = 81; B = "22/10/2011"; C = "MenuManager"; D = "Mapping Item to database : hình như không có thao tác cho phần này"; N = "ThiVT" }
)

foreach ($rd in $rowsData) {
    $r = $rd.Row
    $ws.Range("B$r").Value = $rd.B
    $ws.Range("C$r").Value = $rd.C
    $ws.Range("D$r").Value = $rd.D
    $ws.Range("E$r").Value = "Error"
    $ws.Range("H$r").Value = "Normal"
    $ws.Range("I$r").Value = "Detailed design"
    $ws.Range("K$r").Value = "SDD_MenuManager.xls, revision 116"
    $ws.Range("M$r").Value = "KhoaVT"
    $ws.Range("N$r").Value = $rd.N
    $ws.Range("Q$r").Value = "Sai sót của developer"
    $ws.Range("S$r").Value = "22/10/2011"
    $ws.Rows.Item($r).RowHeight = 30
}

# --- Data validation for the new rows --------------------------------------
# Existing rules on E,F,G,H,I,J,L,R over 61:103 / 63:103 now stop at row 75
# and resume at row 82 (rows 76:81 get their own dedicated rules below).

$ws.Range("E5:E27,E61:E75,E82:E103").Validation.Delete()
$ws.Range("E5:E27,E61:E75,E82:E103").Validation.Add(3, 1, 1, "=$X$7:$X$17")
$ws.Range("E5:E27,E61:E75,E82:E103").Validation.IgnoreBlank = $true
$ws.Range("E5:E27,E61:E75,E82:E103").Validation.InCellDropdown = $true
$ws.Range("E5:E27,E61:E75,E82:E103").Validation.ShowInput = $true
$ws.Range("E5:E27,E61:E75,E82:E103").Validation.ShowError = $true

$ws.Range("F5:F27,F63:F75,F82:F103").Validation.Delete()
$ws.Range("F5:F27,F63:F75,F82:F103").Validation.Add(3, 1, 1, "=$Y$7:$Y$21")
$ws.Range("F5:F27,F63:F75,F82:F103").Validation.IgnoreBlank = $true
$ws.Range("F5:F27,F63:F75,F82:F103").Validation.InCellDropdown = $true
$ws.Range("F5:F27,F63:F75,F82:F103").Validation.ShowInput = $true
$ws.Range("F5:F27,F63:F75,F82:F103").Validation.ShowError = $true

$ws.Range("G5:G27,G63:G75,G82:G103").Validation.Delete()
$ws.Range("G5:G27,G63:G75,G82:G103").Validation.Add(3, 1, 1, "=$Z$7:$Z$15")
$ws.Range("G5:G27,G63:G75,G82:G103").Validation.IgnoreBlank = $true
$ws.Range("G5:G27,G63:G75,G82:G103").Validation.InCellDropdown = $true
$ws.Range("G5:G27,G63:G75,G82:G103").Validation.ShowInput = $true
$ws.Range("G5:G27,G63:G75,G82:G103").Validation.ShowError = $true

$ws.Range("H5:H75,H82:H103").Validation.Delete()
$ws.Range("H5:H75,H82:H103").Validation.Add(3, 1, 1, "=$AA$7:$AA$10")
$ws.Range("H5:H75,H82:H103").Validation.IgnoreBlank = $true
$ws.Range("H5:H75,H82:H103").Validation.InCellDropdown = $true
$ws.Range("H5:H75,H82:H103").Validation.ShowInput = $true
$ws.Range("H5:H75,H82:H103").Validation.ShowError = $true

$ws.Range("I5:I75,I82:I103").Validation.Delete()
$ws.Range("I5:I75,I82:I103").Validation.Add(3, 1, 1, "=$AB$7:$AB$26")
$ws.Range("I5:I75,I82:I103").Validation.IgnoreBlank = $true
$ws.Range("I5:I75,I82:I103").Validation.InCellDropdown = $true
$ws.Range("I5:I75,I82:I103").Validation.ShowInput = $true
$ws.Range("I5:I75,I82:I103").Validation.ShowError = $true

$ws.Range("J5:J27,J63:J75,J82:J103").Validation.Delete()
$ws.Range("J5:J27,J63:J75,J82:J103").Validation.Add(3, 1, 1, "=$AC$7:$AC$10")
$ws.Range("J5:J27,J63:J75,J82:J103").Validation.IgnoreBlank = $true
$ws.Range("J5:J27,J63:J75,J82:J103").Validation.InCellDropdown = $true
$ws.Range("J5:J27,J63:J75,J82:J103").Validation.ShowInput = $true
$ws.Range("J5:J27,J63:J75,J82:J103").Validation.ShowError = $true

$ws.Range("L5:L27,L63:L75,L82:L103").Validation.Delete()
$ws.Range("L5:L27,L63:L75,L82:L103").Validation.Add(3, 1, 1, "=$AD$7:$AD$19")
$ws.Range("L5:L27,L63:L75,L82:L103").Validation.IgnoreBlank = $true
$ws.Range("L5:L27,L63:L75,L82:L103").Validation.InCellDropdown = $true
$ws.Range("L5:L27,L63:L75,L82:L103").Validation.ShowInput = $true
$ws.Range("L5:L27,L63:L75,L82:L103").Validation.ShowError = $true

$ws.Range("R5:R27,R63:R75,R82:R103").Validation.Delete()
$ws.Range("R5:R27,R63:R75,R82:R103").Validation.Add(3, 1, 1, "=$AF$7:$AF$13")
$ws.Range("R5:R27,R63:R75,R82:R103").Validation.IgnoreBlank = $true
$ws.Range("R5:R27,R63:R75,R82:R103").Validation.InCellDropdown = $true
$ws.Range("R5:R27,R63:R75,R82:R103").Validation.ShowInput = $true
$ws.Range("R5:R27,R63:R75,R82:R103").Validation.ShowError = $true

# New validation rules introduced for rows 76:81
$ws.Range("H76:H81").Validation.Add(3, 1, 1, "=$AA$5:$AA$5")
$ws.Range("H76:H81").Validation.IgnoreBlank = $true
$ws.Range("H76:H81").Validation.InCellDropdown = $true
$ws.Range("H76:H81").Validation.ShowInput = $true
$ws.Range("H76:H81").Validation.ShowError = $true

$ws.Range("I76:I81").Validation.Add(3, 1, 1, "=$AB$5:$AB$20")
$ws.Range("I76:I81").Validation.IgnoreBlank = $true
$ws.Range("I76:I81").Validation.InCellDropdown = $true
$ws.Range("I76:I81").Validation.ShowInput = $true
$ws.Range("I76:I81").Validation.ShowError = $true

$ws.Range("J76:J78,L76:L78,R76:R78,E76:G78").Validation.Add(3, 1, 1, "=#REF!")
$ws.Range("J76:J78,L76:L78,R76:R78,E76:G78").Validation.IgnoreBlank = $true
$ws.Range("J76:J78,L76:L78,R76:R78,E76:G78").Validation.InCellDropdown = $true
$ws.Range("J76:J78,L76:L78,R76:R78,E76:G78").Validation.ShowInput = $true
$ws.Range("J76:J78,L76:L78,R76:R78,E76:G78").Validation.ShowError = $true

$ws.Range("J79:J81,L79:L81,R79:R80,E79:G80,E81").Validation.Add(3, 1, 1, "=#REF!")
$ws.Range("J79:J81,L79:L81,R79:R80,E79:G80,E81").Validation.IgnoreBlank = $true
$ws.Range("J79:J81,L79:L81,R79:R80,E79:G80,E81").Validation.InCellDropdown = $true
$ws.Range("J79:J81,L79:L81,R79:R80,E79:G80,E81").Validation.ShowInput = $true
$ws.Range("J79:J81,L79:L81,R79:R80,E79:G80,E81").Validation.ShowError = $true

$ws.Range("R81,F81:G81").Validation.Add(3, 1, 1, "=#REF!")
$ws.Range("R81,F81:G81").Validation.IgnoreBlank = $true
$ws.Range("R81,F81:G81").Validation.InCellDropdown = $true
$ws.Range("R81,F81:G81").Validation.ShowInput = $true
$ws.Range("R81,F81:G81").Validation.ShowError = $true

# --- Restore window view (scroll position / selection) ---------------------
$ws.Range("Q81").Select()
